# feat: add 2022-Q3 data
#
# 1. Update the "总计" (overview) sheet: insert a new "2022-Q3" entry as the
#    first data row, pushing every other quarter down one row and adding
#    2020-Q4 as a brand-new trailing row.
# 2. Insert a brand-new "2022-Q3" worksheet (cloned from the existing
#    "2022-Q2" sheet so headers / styles match) right after the "总计" sheet,
#    populated with the new quarter's fund holdings.
# All the other quarter sheets (2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2,
# 2020-Q4) keep their names & data untouched - they simply shift one tab to
# the right because of the newly inserted sheet.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $addr, $val) {
    # Force the cell to be stored as TEXT (even when the value looks like a
    # number, e.g. fund codes "900090" / "003567" or numeric-looking
    # percentages "51.50") by using Excel's leading-apostrophe text prefix,
    # then resetting the cell style back to Normal so no stray "quote
    # prefix" style/formatting is left behind.
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. "总计" overview sheet (first sheet) - rewrite rows 2-8
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item(1)

$overviewRows = @(
    @(0, "2022-Q3", 4,  5.31),
    @(1, "2022-Q2", 4,  12.09),
    @(2, "2022-Q1", 7,  15.86),
    @(3, "2021-Q4", 10, 18.55),
    @(4, "2021-Q3", 26, 13.2),
    @(5, "2021-Q2", 2,  0.24),
    @(6, "2020-Q4", 5,  1.42)
)

# Row 8 is brand new (sheet used to stop at row 7) - clone the style from an
# existing "A" column cell (s="2") onto A8 before writing its value so the
# new row matches the formatting of every other data row.
$overview.Range("A7").Copy() | Out-Null
$overview.Range("A8").PasteSpecial(-4122) | Out-Null

$r = 2
foreach ($row in $overviewRows) {
    $overview.Cells.Item($r, 1).Value = $row[0]
    $overview.Cells.Item($r, 2).Value = $row[1]
    $overview.Cells.Item($r, 3).Value = $row[2]
    $overview.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Brand-new "2022-Q3" sheet, cloned from "2022-Q2" (position 2) so it
#    inherits identical headers/styles, inserted right before it.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($q2, $null) | Out-Null
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

Set-TextCell $q3 "B2" "900090"
Set-TextCell $q3 "C2" "中信卓越成长两年持有期混合B"
Set-TextCell $q3 "D2" "51.50"
Set-TextCell $q3 "E2" "93.14"
Set-TextCell $q3 "F2" "4.47"
Set-TextCell $q3 "G2" "2.3020"
$q3.Range("H2").Value = 5

Set-TextCell $q3 "B3" "003567"
Set-TextCell $q3 "C3" "华夏行业景气混合"
Set-TextCell $q3 "D3" "115.66"
Set-TextCell $q3 "E3" "88.33"
Set-TextCell $q3 "F3" "1.87"
Set-TextCell $q3 "G3" "2.1628"
$q3.Range("H3").Value = 10

Set-TextCell $q3 "B4" "900010"
Set-TextCell $q3 "C4" "中信卓越成长两年持有期混合A"
Set-TextCell $q3 "D4" "14.24"
Set-TextCell $q3 "E4" "93.14"
Set-TextCell $q3 "F4" "4.47"
Set-TextCell $q3 "G4" "0.6365"
$q3.Range("H4").Value = 5

Set-TextCell $q3 "B5" "900100"
Set-TextCell $q3 "C5" "中信卓越成长两年持有期混合C"
Set-TextCell $q3 "D5" "4.61"
Set-TextCell $q3 "E5" "93.14"
Set-TextCell $q3 "F5" "4.47"
Set-TextCell $q3 "G5" "0.2061"
$q3.Range("H5").Value = 5

$overview.Select()
$overview.Range("A1").Select()
